# This commit adds one new weekly observation for "Zanahoria" at the
# "Feria Lagunitas de Puerto Montt" market. The new record belongs right
# after the current row 586 (the data is kept in reverse-chronological
# order), so a new row is inserted at row 587 and every existing record
# that used to live at rows 587-615 is pushed down to rows 588-616.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 587, shifting rows 587:615 down to 588:616
# (this also carries the existing formatting, e.g. the date number
# format on column D, down along with the shifted rows).
$ws.Rows("587:587").Insert()

# Populate the newly-inserted row with the new observation.
$ws.Cells.Item(587, 1).Value  = 4
$ws.Cells.Item(587, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(587, 3).Value  = "Los Lagos"
$ws.Cells.Item(587, 4).Value  = 45147
$ws.Cells.Item(587, 5).Value  = 10
$ws.Cells.Item(587, 6).Value  = 100114013
$ws.Cells.Item(587, 7).Value  = "Zanahoria"
$ws.Cells.Item(587, 8).Value  = "Sin especificar"
$ws.Cells.Item(587, 9).Value  = "Primera"
$ws.Cells.Item(587, 10).Value = 150
$ws.Cells.Item(587, 11).Value = 7500
$ws.Cells.Item(587, 12).Value = 7500
$ws.Cells.Item(587, 13).Value = 7500
$ws.Cells.Item(587, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(587, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(587, 16).Value = 375
$ws.Cells.Item(587, 17).Value = 20
$ws.Cells.Item(587, 18).Value = "Hortaliza"
